# Fruta / hortaliza, semanal
#
# The source data table (Hortaliza, Femacal de La Calera - Pepino ensalada)
# gets a new weekly observation. A new row is inserted at row 191, pushing
# the existing rows 191-225 down to 192-226, and the new row is populated
# with the latest weekly values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 191 (shifts rows 191:225 down to 192:226,
# extending the used range / dimension to A1:R226, same as Excel's
# Rows("191:191").Insert or Range("A191").EntireRow.Insert()).
$ws.Range("A191").EntireRow.Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A191").Value = 3
$ws.Range("B191").Value = "Femacal de La Calera"
$ws.Range("C191").Value = "Coquimbo"
$ws.Range("D191").Value = 44504
$ws.Range("E191").Value = 5
$ws.Range("F191").Value = 100112043
$ws.Range("G191").Value = "Pepino ensalada"
$ws.Range("H191").Value = "Sin especificar"
$ws.Range("I191").Value = "Primera"
$ws.Range("J191").Value = 115
$ws.Range("K191").Value = 7000
$ws.Range("L191").Value = 7500
$ws.Range("M191").Value = 7217
$ws.Range("N191").Value = "`$/caja 70 unidades"
$ws.Range("O191").Value = "Región de Arica y Parinacota"
$ws.Range("P191").Value = 103
$ws.Range("Q191").Value = 70
$ws.Range("R191").Value = "Hortaliza"
